# Remove the "valorDescarga"/"opcionSubmenu" header and their sample data values
# in columns M and N (the feature no longer validates these fields), and select
# the cleared range M1:N2 to reflect it in the sheet view.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M1:N2").ClearContents()
$ws.Range("M1:N2").Select()
